{"js": "// Mapping of old equation text -> new equation text, in document order\n// (row-major over the 20x5 answers table).\nconst PAIRS = [[\"9+13=22\", \"99-57=42\"], [\"41+20=61\", \"17+81=98\"], [\"45+9=54\", \"53-49=4\"], [\"14+79=93\", \"61-10=51\"], [\"30+23=53\", \"54+33=87\"], [\"38-37=1\", \"45-19=26\"], [\"88-16=72\", \"62+8=70\"], [\"32-14=18\", \"55-19=36\"], [\"23+24=47\", \"71+12=83\"], [\"56+36=92\", \"33+16=49\"], [\"84+14=98\", \"56-50=6\"], [\"64+11=75\", \"19+18=37\"], [\"26+72=98\", \"87-38=49\"], [\"69-15=54\", \"91+6=97\"], [\"80-78=2\", \"89-84=5\"], [\"45+49=94\", \"31+56=87\"], [\"8+58=66\", \"84-2=82\"], [\"0+11=11\", \"58-29=29\"], [\"28-0=28\", \"31+61=92\"], [\"10+11=21\", \"37+0=37\"], [\"95-86=9\", \"96-53=43\"], [\"53-7=46\", \"34+12=46\"], [\"47-28=19\", \"62-27=35\"], [\"98-10=88\", \"61+6=67\"], [\"63-29=34\", \"78-75=3\"], [\"46+6=52\", \"30+68=98\"], [\"52-4=48\", \"41-30=11\"], [\"23+49=72\", \"85-10=75\"], [\"79-10=69\", \"25+18=43\"], [\"97-1=96\", \"30+19=49\"], [\"2+3=5\", \"56+18=74\"], [\"93-41=52\", \"32-19=13\"], [\"17+6=23\", \"33-30=3\"], [\"41+46=87\", \"91-33=58\"], [\"46-38=8\", \"87-84=3\"], [\"45+34=79\", \"43+30=73\"], [\"85-41=44\", \"10+46=56\"], [\"26+61=87\", \"11+36=47\"], [\"55-25=30\", \"89-3=86\"], [\"91-75=16\", \"46+20=66\"], [\"12+70=82\", \"54+35=89\"], [\"14+81=95\", \"23-1=22\"], [\"94-33=61\", \"66+17=83\"], [\"51-44=7\", \"55-41=14\"], [\"44+39=83\", \"68-15=53\"], [\"68+16=84\", \"49+34=83\"], [\"57-4=53\", \"54-53=1\"], [\"11+49=60\", \"75-14=61\"], [\"57-36=21\", \"43+51=94\"], [\"24+47=71\", \"56-19=37\"], [\"88+6=94\", \"43-35=8\"], [\"40-16=24\", \"70-22=48\"], [\"5+60=65\", \"6+40=46\"], [\"3+14=17\", \"71-22=49\"], [\"44+28=72\", \"80-62=18\"], [\"13+61=74\", \"94-24=70\"], [\"8+2=10\", \"84-56=28\"], [\"57+0=57\", \"88-77=11\"], [\"71-39=32\", \"33+6=39\"], [\"64-62=2\", \"77-18=59\"], [\"56+0=56\", \"44-8=36\"], [\"56+29=85\", \"38+45=83\"], [\"21+37=58\", \"27-13=14\"], [\"37-25=12\", \"31+60=91\"], [\"95-32=63\", \"68-49=19\"], [\"70-66=4\", \"28+52=80\"], [\"52-13=39\", \"46+13=59\"], [\"45+39=84\", \"73+14=87\"], [\"82-62=20\", \"46-25=21\"], [\"62-23=39\", \"64-40=24\"], [\"21+31=52\", \"83-34=49\"], [\"91-48=43\", \"93-59=34\"], [\"29+56=85\", \"90-48=42\"], [\"87-6=81\", \"21+43=64\"], [\"45-8=37\", \"29+47=76\"], [\"21+48=69\", \"1+96=97\"], [\"61+11=72\", \"84-70=14\"], [\"47+46=93\", \"95-75=20\"], [\"77-72=5\", \"82-58=24\"], [\"9+1=10\", \"59-38=21\"], [\"66+1=67\", \"56+39=95\"], [\"2+10=12\", \"73-64=9\"], [\"7+62=69\", \"17+77=94\"], [\"25+41=66\", \"70-14=56\"], [\"49+40=89\", \"57-46=11\"], [\"16+78=94\", \"2+79=81\"], [\"85-16=69\", \"83+7=90\"], [\"56+20=76\", \"42+53=95\"], [\"60-55=5\", \"69-49=20\"], [\"13+16=29\", \"99-54=45\"], [\"65-37=28\", \"74-30=44\"], [\"16+18=34\", \"60+36=96\"], [\"83-22=61\", \"1+92=93\"], [\"95-65=30\", \"12+71=83\"], [\"48-13=35\", \"1+60=61\"], [\"93-88=5\", \"35+7=42\"], [\"80-10=70\", \"95-87=8\"], [\"24+58=82\", \"18+11=29\"], [\"64-8=56\", \"26+9=35\"], [\"90-40=50\", \"37+17=54\"]];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document\");\n}\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows;\n// Load cells for every row first.\nfor (let r = 0; r < rows.items.length; r++) {\n  rows.items[r].cells.load(\"items\");\n}\nawait context.sync();\n\n// Flatten all cells in row-major order and load their current values.\nconst cells = [];\nfor (let r = 0; r < rows.items.length; r++) {\n  const rowCells = rows.items[r].cells;\n  for (let c = 0; c < rowCells.items.length; c++) {\n    cells.push(rowCells.items[c]);\n  }\n}\nfor (const cell of cells) {\n  cell.load(\"value\");\n}\nawait context.sync();\n\nif (cells.length !== PAIRS.length) {\n  throw new Error(\n    `Expected ${PAIRS.length} answer cells, found ${cells.length}`\n  );\n}\n\n// Replace each cell's text with its mapped new value. Cells are walked in\n// row-major document order (matching PAIRS) and the current value is\n// checked against the expected \"old\" value first, so the script fails\n// loudly instead of silently mis-mapping if the table ever changes shape.\nfor (let i = 0; i < cells.length; i++) {\n  const cell = cells[i];\n  const [oldVal, newVal] = PAIRS[i];\n  const current = (cell.value || \"\").trim();\n  if (current !== oldVal) {\n    throw new Error(\n      `Cell ${i}: expected \"${oldVal}\" but found \"${current}\"`\n    );\n  }\n  cell.getRange().insertText(newVal, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Old -> New equation text for every cell of the single 20x5 answers table,\n# in row-major order (matches $d.Tables.Item(1).Cell(r, c) iteration below).\n$Old = @(\n  \"9+13=22\",\n  \"41+20=61\",\n  \"45+9=54\",\n  \"14+79=93\",\n  \"30+23=53\",\n  \"38-37=1\",\n  \"88-16=72\",\n  \"32-14=18\",\n  \"23+24=47\",\n  \"56+36=92\",\n  \"84+14=98\",\n  \"64+11=75\",\n  \"26+72=98\",\n  \"69-15=54\",\n  \"80-78=2\",\n  \"45+49=94\",\n  \"8+58=66\",\n  \"0+11=11\",\n  \"28-0=28\",\n  \"10+11=21\",\n  \"95-86=9\",\n  \"53-7=46\",\n  \"47-28=19\",\n  \"98-10=88\",\n  \"63-29=34\",\n  \"46+6=52\",\n  \"52-4=48\",\n  \"23+49=72\",\n  \"79-10=69\",\n  \"97-1=96\",\n  \"2+3=5\",\n  \"93-41=52\",\n  \"17+6=23\",\n  \"41+46=87\",\n  \"46-38=8\",\n  \"45+34=79\",\n  \"85-41=44\",\n  \"26+61=87\",\n  \"55-25=30\",\n  \"91-75=16\",\n  \"12+70=82\",\n  \"14+81=95\",\n  \"94-33=61\",\n  \"51-44=7\",\n  \"44+39=83\",\n  \"68+16=84\",\n  \"57-4=53\",\n  \"11+49=60\",\n  \"57-36=21\",\n  \"24+47=71\",\n  \"88+6=94\",\n  \"40-16=24\",\n  \"5+60=65\",\n  \"3+14=17\",\n  \"44+28=72\",\n  \"13+61=74\",\n  \"8+2=10\",\n  \"57+0=57\",\n  \"71-39=32\",\n  \"64-62=2\",\n  \"56+0=56\",\n  \"56+29=85\",\n  \"21+37=58\",\n  \"37-25=12\",\n  \"95-32=63\",\n  \"70-66=4\",\n  \"52-13=39\",\n  \"45+39=84\",\n  \"82-62=20\",\n  \"62-23=39\",\n  \"21+31=52\",\n  \"91-48=43\",\n  \"29+56=85\",\n  \"87-6=81\",\n  \"45-8=37\",\n  \"21+48=69\",\n  \"61+11=72\",\n  \"47+46=93\",\n  \"77-72=5\",\n  \"9+1=10\",\n  \"66+1=67\",\n  \"2+10=12\",\n  \"7+62=69\",\n  \"25+41=66\",\n  \"49+40=89\",\n  \"16+78=94\",\n  \"85-16=69\",\n  \"56+20=76\",\n  \"60-55=5\",\n  \"13+16=29\",\n  \"65-37=28\",\n  \"16+18=34\",\n  \"83-22=61\",\n  \"95-65=30\",\n  \"48-13=35\",\n  \"93-88=5\",\n  \"80-10=70\",\n  \"24+58=82\",\n  \"64-8=56\",\n  \"90-40=50\"\n)\n\n$New = @(\n  \"99-57=42\",\n  \"17+81=98\",\n  \"53-49=4\",\n  \"61-10=51\",\n  \"54+33=87\",\n  \"45-19=26\",\n  \"62+8=70\",\n  \"55-19=36\",\n  \"71+12=83\",\n  \"33+16=49\",\n  \"56-50=6\",\n  \"19+18=37\",\n  \"87-38=49\",\n  \"91+6=97\",\n  \"89-84=5\",\n  \"31+56=87\",\n  \"84-2=82\",\n  \"58-29=29\",\n  \"31+61=92\",\n  \"37+0=37\",\n  \"96-53=43\",\n  \"34+12=46\",\n  \"62-27=35\",\n  \"61+6=67\",\n  \"78-75=3\",\n  \"30+68=98\",\n  \"41-30=11\",\n  \"85-10=75\",\n  \"25+18=43\",\n  \"30+19=49\",\n  \"56+18=74\",\n  \"32-19=13\",\n  \"33-30=3\",\n  \"91-33=58\",\n  \"87-84=3\",\n  \"43+30=73\",\n  \"10+46=56\",\n  \"11+36=47\",\n  \"89-3=86\",\n  \"46+20=66\",\n  \"54+35=89\",\n  \"23-1=22\",\n  \"66+17=83\",\n  \"55-41=14\",\n  \"68-15=53\",\n  \"49+34=83\",\n  \"54-53=1\",\n  \"75-14=61\",\n  \"43+51=94\",\n  \"56-19=37\",\n  \"43-35=8\",\n  \"70-22=48\",\n  \"6+40=46\",\n  \"71-22=49\",\n  \"80-62=18\",\n  \"94-24=70\",\n  \"84-56=28\",\n  \"88-77=11\",\n  \"33+6=39\",\n  \"77-18=59\",\n  \"44-8=36\",\n  \"38+45=83\",\n  \"27-13=14\",\n  \"31+60=91\",\n  \"68-49=19\",\n  \"28+52=80\",\n  \"46+13=59\",\n  \"73+14=87\",\n  \"46-25=21\",\n  \"64-40=24\",\n  \"83-34=49\",\n  \"93-59=34\",\n  \"90-48=42\",\n  \"21+43=64\",\n  \"29+47=76\",\n  \"1+96=97\",\n  \"84-70=14\",\n  \"95-75=20\",\n  \"82-58=24\",\n  \"59-38=21\",\n  \"56+39=95\",\n  \"73-64=9\",\n  \"17+77=94\",\n  \"70-14=56\",\n  \"57-46=11\",\n  \"2+79=81\",\n  \"83+7=90\",\n  \"42+53=95\",\n  \"69-49=20\",\n  \"99-54=45\",\n  \"74-30=44\",\n  \"60+36=96\",\n  \"1+92=93\",\n  \"12+71=83\",\n  \"1+60=61\",\n  \"35+7=42\",\n  \"95-87=8\",\n  \"18+11=29\",\n  \"26+9=35\",\n  \"37+17=54\"\n)\n\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\nif (($rowCount * $colCount) -ne $Old.Count) {\n    throw \"Expected $($Old.Count) cells but table has $rowCount x $colCount = $($rowCount * $colCount)\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $table.Cell($r, $c)\n        $current = $cell.Range.Text.TrimEnd([char]7, [char]13)\n        $expected = $Old[$i]\n        if ($current -ne $expected) {\n            throw \"Cell ($r,$c): expected '$expected' but found '$current'\"\n        }\n        $cell.Range.Text = $New[$i]\n        $i++\n    }\n}\n\nWrite-Output \"Replaced $i cells\"\n"}
